# Generate Report for Handoff
#
# For the six "Ready for handoff" rows (source files 26c0e6ef, 28003d56,
# 6fb53695, 8a19dfb0, 9d52b4be, 9e0095ee -> sheet rows 7,8,9,11,12,13):
#   - Overview sheet: "Latest HO Xliff Generate Date" (col G) advances
#     from 2016-08-26 16:22:14 to 2016-08-26 16:22:34
#   - zh-cn sheet: "Latest Handoff Datetime" (col H) advances from
#     2016-08-26 16:22:08 to 2016-08-26 16:22:29, and "Priority" (col E)
#     is set to "ht"
#   - de-de sheet: "Latest Handoff Datetime" (col H) advances from
#     2016-08-26 16:22:14 to 2016-08-26 16:22:34, and "Priority" (col E)
#     is set to "ht"

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$handoffRows = @(7, 8, 9, 11, 12, 13)

foreach ($row in $handoffRows) {
    # Overview!G<row> - Latest HO Xliff Generate Date
    $wsOverview.Cells.Item($row, 7).Value2 = "2016-08-26 16:22:34"

    # zh-cn!H<row> - Latest Handoff Datetime
    $wsZhCn.Cells.Item($row, 8).Value2 = "2016-08-26 16:22:29"
    # zh-cn!E<row> - Priority
    $wsZhCn.Cells.Item($row, 5).Value2 = "ht"

    # de-de!H<row> - Latest Handoff Datetime
    $wsDeDe.Cells.Item($row, 8).Value2 = "2016-08-26 16:22:34"
    # de-de!E<row> - Priority
    $wsDeDe.Cells.Item($row, 5).Value2 = "ht"
}
